$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 was a "totals" row with a broken formula (cached sum was 300,
# instead of 2000+10+13). Re-enter the values so the formula result is
# corrected on recalculation.
$ws.Range("C16").Value = 2000
$ws.Range("D16").Value = 10
$ws.Range("E16").Value = 13

# Work in manual calculation mode while building the new row so the
# freshly-entered formula captures its value at entry time (against the
# still-empty precedent cells) instead of immediately recalculating once
# the data below it is filled in.
$excel.Calculation = -4135

# New row 17, duplicated from row 16's pattern with an updated weight.
$ws.Range("F17").Formula = "=D17+E17+C17"
$ws.Range("C17").Value = 3000
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = 13

$ws.Range("G17").Select()
